$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Templates")

# Insert a new column before the old column F ("Description"), shifting
# Description from F -> G and leaving a blank column F for the new
# "File Name" data.
$ws.Columns("F:F").Insert()

# Copy the formatting (styles) of column A (Template Number) onto the new
# column F for rows 1-40, matching header style on row 1 and the
# centered-text style on the data rows.
$ws.Range("A1:A40").Copy()
$ws.Range("F1:F40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give column F the same width as the other narrow text columns (15
# characters, same as column E).
$ws.Columns("F:F").ColumnWidth = 14.17

# Header for the new column.
$ws.Range("F1").Value = "File Name"

# File Name values - one per template row, mirroring the Template Number
# in column A for every row except row 32, where template "031" reuses
# the "007" template file for a new/second position.
$fileNames = @{
    2  = "001"
    3  = "002"
    4  = "003"
    5  = "004"
    6  = "005"
    7  = "006"
    8  = "007"
    9  = "008"
    10 = "009"
    11 = "010"
    12 = "011"
    13 = "012"
    14 = "013"
    15 = "014"
    16 = "015"
    17 = "016"
    18 = "017"
    19 = "018"
    20 = "019"
    21 = "020"
    22 = "021"
    23 = "022"
    24 = "023"
    25 = "024"
    26 = "025"
    27 = "026"
    28 = "027"
    29 = "028"
    30 = "029"
    31 = "030"
    32 = "007"
    33 = "032"
    34 = "033"
    35 = "034"
    36 = "035"
    37 = "036"
    38 = "037"
    39 = "038"
    40 = "039"
}

foreach ($row in $fileNames.Keys) {
    $ws.Cells.Item($row, 6).Value = $fileNames[$row]
}

# Restore the sheet view: scroll back to the top and select the new
# bottom-right corner of the data (G40, under the Description column).
$ws.Activate()
$ws.Range("G40").Select()

Write-Output "done"
